$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.594.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.68%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.880.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.13%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.95%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'316.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.26%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.71%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.5091"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3905"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.22%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.08425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.11%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.34%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'41.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'6.235"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.49%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.874.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.93%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'20.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.40%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'7.252"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.86%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'1.011"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.11%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +1.06%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'91.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.37%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06738"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.29%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.86%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +0.77%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.942"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.35%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'28.641.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.72%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.66%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.244"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.16%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.088.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.99%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'162.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.57%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.64%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.368"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.66%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'126.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.18%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.1047"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.60%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +1.16%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.26%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.624"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.84%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +1.66%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.06556"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.57%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.2166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.18%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'8.866"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.96%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'5.086"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.06%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.198"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.83%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.258"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.95%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +0.29%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'11.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.93%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +0.69%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.6049"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.87%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +0.91%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.701"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.99%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.014"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.80%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.219"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.52%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'122.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.08%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -10.69%  "
$ws.Range("E51").Style = "Normal"
